$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 12
$ws.Range("H2").Value = 1.28
$ws.Range("I2").Value = 1.34
$ws.Range("K2").Value = 6.6
$ws.Range("P2").Value = 2.36
$ws.Range("Q2").Value = 1.62
$ws.Range("BH2").Value = "2026-02-24 12:21:44"

# Row 3
$ws.Range("H3").Value = 4.2
$ws.Range("I3").Value = 4.8
$ws.Range("J3").Value = 3.55
$ws.Range("K3").Value = 4
$ws.Range("P3").Value = 1.94
$ws.Range("Q3").Value = 1.9
$ws.Range("BH3").Value = "2026-02-24 12:21:44"

# Row 4
$ws.Range("BH4").Value = "2026-02-24 12:21:44"

# Row 5
$ws.Range("F5").Value = 1.83
$ws.Range("Q5").Value = 2.42
$ws.Range("BH5").Value = "2026-02-24 12:21:44"

# Row 6
$ws.Range("BH6").Value = "2026-02-24 12:21:44"

# Row 7
$ws.Range("K7").Value = 3.6
$ws.Range("BH7").Value = "2026-02-24 12:21:44"

# Row 8
$ws.Range("BH8").Value = "2026-02-24 12:21:44"

# Row 9
$ws.Range("F9").Value = 1.49
$ws.Range("P9").Value = 1.66
$ws.Range("BH9").Value = "2026-02-24 12:21:44"

# Row 10
$ws.Range("BH10").Value = "2026-02-24 12:21:44"
